$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Wishlist Maia & Clara"

# Mark rows as reserved/bought (column E = "Y") and hide them
$ws.Range("E23").Value = "Y"
$ws.Range("E28:E33").Value = "Y"

$ws.Rows("23:26").Hidden = $true
$ws.Rows("28:33").Hidden = $true

# Reapply the autofilter range to cover the new extent, restoring the
# "blank" filter criteria on column E (5th column of the range)
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:E34").AutoFilter(5, @(""), 7)

# Keep the _FilterDatabase defined name range in sync with the autofilter
$fdb = $wb.Names.Item("_xlnm._FilterDatabase")
$fdb.RefersTo = "='Wishlist Maia & Clara'!`$A`$1:`$E`$34"

[void]$ws.Range("E58").Select()
